# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that any contributor name containing "system" (case-insensitive, covers
# both "system" and "System") is moved to the end of the comma-separated
# list, while the relative order of all entries is otherwise preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) { continue }
    if ($text.IndexOf(",") -lt 0) { continue }   # only one contributor, nothing to reorder

    $parts = $text -split ", "

    $nonSystem = @()
    $systemParts = @()
    foreach ($p in $parts) {
        if ($p.ToLower().Contains("system")) {
            $systemParts += $p
        } else {
            $nonSystem += $p
        }
    }

    $newVal = ($nonSystem + $systemParts) -join ", "

    if ($newVal -ne $text) {
        $cell.Value = $newVal
    }
}
